$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'" + '29.488.99'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.52%  '

$ws.Range('D3').Value = "'" + '1.915.13'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.53%  '

$ws.Range('D4').Value = "'" + '1.007'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.62%  '

$ws.Range('D5').Value = "'" + '325.47'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.87%  '

$ws.Range('D6').Value = "'" + '1.005'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.44%  '

$ws.Range('D7').Value = "'" + '0.4845'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +2.64%  '

$ws.Range('D8').Value = "'" + '0.4073'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +1.08%  '

$ws.Range('D9').Value = "'" + '0.08176'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.90%  '

$ws.Range('D10').Value = "'" + '1.015'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +2.20%  '

$ws.Range('D11').Value = "'" + '23.74'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +5.01%  '

$ws.Range('D12').Value = "'" + '1.909.68'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.75%  '

$ws.Range('D13').Value = "'" + '6.039'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +2.96%  '

$ws.Range('D14').Value = "'" + '7.188'
$ws.Range('D14').Style = 'Normal'

$ws.Range('E15').Value = '  +2.13%  '

$ws.Range('D16').Value = "'" + '0.06773'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +2.45%  '

$ws.Range('D17').Value = "'" + '1.007'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.51%  '

$ws.Range('E18').Value = '  +1.10%  '

$ws.Range('D19').Value = "'" + '17.73'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.99%  '

$ws.Range('E20').Value = '  +0.53%  '

$ws.Range('D21').Value = "'" + '29.501.68'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.57%  '

$ws.Range('E22').Value = '  +2.12%  '

$ws.Range('D23').Value = "'" + '11.75'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +2.53%  '

$ws.Range('D24').Value = "'" + '2.178'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.10%  '

$ws.Range('D25').Value = "'" + '2.129.24'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.15%  '

$ws.Range('B26').Value = 'Monero'
$ws.Range('C26').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D26').Value = "'" + '156.44'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +1.34%  '

$ws.Range('B27').Value = 'InternetComputer(DFINITY)'
$ws.Range('C27').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D27').Value = "'" + '6.524'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +8.04%  '

$ws.Range('D28').Value = "'" + '20.12'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.98%  '

$ws.Range('D29').Value = "'" + '2.126'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +1.54%  '

$ws.Range('D30').Value = "'" + '120.72'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +2.61%  '

$ws.Range('D31').Value = "'" + '1.028'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -3.60%  '

$ws.Range('D32').Value = "'" + '0.09554'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +1.03%  '

$ws.Range('D33').Value = "'" + '5.505'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +2.78%  '

$ws.Range('B34').Value = 'ARBITRUM'
$ws.Range('C34').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D34').Value = "'" + '1.397'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.94%  '

$ws.Range('B35').Value = 'HuobiToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D35').Value = "'" + '3.555'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.04%  '

$ws.Range('D36').Value = "'" + '0.02280'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +1.55%  '

$ws.Range('D37').Value = "'" + '0.06128'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.83%  '

$ws.Range('D38').Value = "'" + '1.188'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +1.24%  '

$ws.Range('D39').Value = "'" + '10.89'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +8.01%  '

$ws.Range('D40').Value = "'" + '0.5973'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +2.82%  '

$ws.Range('D41').Value = "'" + '7.997'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.96%  '

$ws.Range('D42').Value = "'" + '0.1854'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +1.27%  '

$ws.Range('B43').Value = 'RenderToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D43').Value = "'" + '2.414'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -3.94%  '

$ws.Range('B44').Value = 'WEMIXToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D44').Value = "'" + '1.279'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.59%  '

$ws.Range('D45').Value = "'" + '0.07630'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.52%  '

$ws.Range('D46').Value = "'" + '12.42'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +2.12%  '

$ws.Range('D47').Value = "'" + '0.5583'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +1.90%  '

$ws.Range('E48').Value = '  +2.84%  '

$ws.Range('D49').Value = "'" + '116.69'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +2.78%  '

$ws.Range('D50').Value = "'" + '72.67'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +2.22%  '

$ws.Range('D51').Value = "'" + '2.413'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +3.05%  '
